$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 761
$ws.Range("G2").Value = 50
$ws.Range("F3").Value = 149
$ws.Range("F5").Value = 394
$ws.Range("F7").Value = 4202
$ws.Range("F8").Value = 344
$ws.Range("F9").Value = 219
$ws.Range("F10").Value = 835
$ws.Range("F11").Value = 790
$ws.Range("F12").Value = 62
$ws.Range("F13").Value = 7
$ws.Range("F14").Value = 525
$ws.Range("F17").Value = 1515
$ws.Range("F18").Value = 1395
$ws.Range("F19").Value = 596
$ws.Range("F21").Value = 164
$ws.Range("F22").Value = 211
$ws.Range("F23").Value = 424
$ws.Range("F27").Value = 521
$ws.Range("F28").Value = 851
$ws.Range("F29").Value = 105
$ws.Range("F31").Value = 144
$ws.Range("F36").Value = 231
$ws.Range("F37").Value = 490
$ws.Range("F38").Value = 55

$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 118
$ws.Range("F6").Value = 87
$ws.Range("F7").Value = 4

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 761
$ws.Range("G3").Value = 50
$ws.Range("F4").Value = 149
$ws.Range("F6").Value = 394
$ws.Range("F9").Value = 4202
$ws.Range("F10").Value = 344
$ws.Range("F11").Value = 219
$ws.Range("F13").Value = 118
$ws.Range("F14").Value = 835
$ws.Range("F15").Value = 790
$ws.Range("F17").Value = 87
$ws.Range("F18").Value = 62
$ws.Range("F19").Value = 7
$ws.Range("F20").Value = 525
$ws.Range("F22").Value = 4
$ws.Range("F24").Value = 1515
$ws.Range("F25").Value = 1395
$ws.Range("F26").Value = 596
$ws.Range("F28").Value = 164
$ws.Range("F29").Value = 211
$ws.Range("F31").Value = 424
$ws.Range("F35").Value = 521
$ws.Range("F36").Value = 851
$ws.Range("F37").Value = 105
$ws.Range("F39").Value = 144
$ws.Range("F44").Value = 231
$ws.Range("F45").Value = 490
$ws.Range("F46").Value = 55
